$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 115 (last data row). Two new match rows
# (116, 117) are appended, continuing the existing table layout/style.
# Copy formatting from the last existing data row down onto the two new
# rows so they inherit the same cell styles (bordered/bold Indice column,
# date-time format on data_partida, etc.) as the rest of the table.
$ws.Range("A115:V115").Copy() | Out-Null
$ws.Range("A116:V117").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 116: Piast Gliwice 0 x 0 Korona Kielce
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = "poland"
$ws.Cells.Item(116, 3).Value = "ekstraklasa"
$ws.Cells.Item(116, 4).Value = "2023-2024"
$ws.Cells.Item(116, 5).Value = 45233.75
$ws.Cells.Item(116, 6).Value = "Piast Gliwice"
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = "Korona Kielce"
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 1.57
$ws.Cells.Item(116, 11).Value = "28/10/2023 12:42"
$ws.Cells.Item(116, 12).Value = 1.58
$ws.Cells.Item(116, 13).Value = "03/11/2023 17:58"
$ws.Cells.Item(116, 14).Value = 3.89
$ws.Cells.Item(116, 15).Value = "28/10/2023 12:42"
$ws.Cells.Item(116, 16).Value = 3.85
$ws.Cells.Item(116, 17).Value = "03/11/2023 17:58"
$ws.Cells.Item(116, 18).Value = 5.68
$ws.Cells.Item(116, 19).Value = "28/10/2023 12:42"
$ws.Cells.Item(116, 20).Value = 6.92
$ws.Cells.Item(116, 21).Value = "03/11/2023 17:58"
$ws.Cells.Item(116, 22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/piast-gliwice-korona-kielce/2LhSMI0N/"

# Row 117: Jagiellonia 4 x 0 Stal Mielec
$ws.Cells.Item(117, 1).Value = 116
$ws.Cells.Item(117, 2).Value = "poland"
$ws.Cells.Item(117, 3).Value = "ekstraklasa"
$ws.Cells.Item(117, 4).Value = "2023-2024"
$ws.Cells.Item(117, 5).Value = 45233.85416666666
$ws.Cells.Item(117, 6).Value = "Jagiellonia"
$ws.Cells.Item(117, 7).Value = 4
$ws.Cells.Item(117, 8).Value = "Stal Mielec"
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 1.49
$ws.Cells.Item(117, 11).Value = "30/10/2023 03:12"
$ws.Cells.Item(117, 12).Value = 1.54
$ws.Cells.Item(117, 13).Value = "03/11/2023 19:55"
$ws.Cells.Item(117, 14).Value = 4.42
$ws.Cells.Item(117, 15).Value = "30/10/2023 03:12"
$ws.Cells.Item(117, 16).Value = 4.37
$ws.Cells.Item(117, 17).Value = "03/11/2023 19:55"
$ws.Cells.Item(117, 18).Value = 6.63
$ws.Cells.Item(117, 19).Value = "30/10/2023 03:12"
$ws.Cells.Item(117, 20).Value = 6.3
$ws.Cells.Item(117, 21).Value = "03/11/2023 19:55"
$ws.Cells.Item(117, 22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/jagiellonia-stal-mielec/U7fKOvWA/"

Write-Host "Appended rows 116-117 to sheet '$($ws.Name)'. New dimension: $($ws.UsedRange.Address())"
